# Uren Registratie 23-02-2016 & story plaatjes toegevoegd presentatie
#
# Row 28 ("Story" row, under the "Sprint 1" block) previously had no hours
# logged. Four hours are now logged for every team member (columns B-G),
# and the cells that already carried the team's "hours-registered" fill
# colour (columns C-G) get the same blue fill that the rest of the filled
# rows in the sheet use. All of the totals/averages above (row 2 block)
# recalculate automatically because they are driven by formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log 4 hours for every person on row 28.
$ws.Range("B28:G28").Value = 4

# Re-apply the "hours logged" blue fill to the newly filled cells
# (same colour already used elsewhere in the sheet: RGB(0,112,192)).
$ws.Range("C28:G28").Interior.Color = 12611584

# Move the active cell/selection the way the author left the workbook.
$ws.Range("J21").Select() | Out-Null
